$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D32").Value = "Individual conditional expectation (ICE)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/373"

$ws.Range("D46").Value = "고나트륨혈증, 저나트륨혈증"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/473"
